$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = -0.1119300937248394
$ws.Range("J18").Value = 0.2456466839623254
$ws.Range("K18").Value = -0.1838583654545527
$ws.Range("L18").Value = 2.610687095238301

$ws.Range("I19").Value = 0.2190844768866997
$ws.Range("J19").Value = 0.5198204364300951
$ws.Range("K19").Value = 0.1047479672860789
$ws.Range("L19").Value = 2.06667398167574
